$p = $ppt.ActivePresentation

# Slide 3: "EINLEITUNG INS PROJEKT" -> "Einleitung ins Projekt" (title shape)
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "EINLEITUNG INS PROJEKT") {
            $shp.TextFrame.TextRange.Text = "Einleitung ins Projekt"
        }
    }
}

# Slide 18: "Gerne beantworte wir jetzt Ihre Fragen." -> "Gerne beantworten wir jetzt Ihre Fragen."
$s18 = $p.Slides.Item(18)
for ($i = 1; $i -le $s18.Shapes.Count; $i++) {
    $shp = $s18.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Gerne beantworte wir jetzt Ihre Fragen.") {
            $shp.TextFrame.TextRange.Text = "Gerne beantworten wir jetzt Ihre Fragen."
        }
    }
}
